$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 300
$ws1.Range("F4").Value = 219
$ws1.Range("F5").Value = 246
$ws1.Range("F7").Value = 7367
$ws1.Range("F9").Value = 67
$ws1.Range("F10").Value = 3154
$ws1.Range("F12").Value = 577
$ws1.Range("F13").Value = 587
$ws1.Range("F17").Value = 746
$ws1.Range("F20").Value = 198
$ws1.Range("F26").Value = 1069
$ws1.Range("F28").Value = 115
$ws1.Range("F29").Value = 2101
$ws1.Range("F30").Value = 616
$ws1.Range("F31").Value = 24

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 304

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 404

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 404
$ws4.Range("F3").Value = 300
$ws4.Range("F5").Value = 219
$ws4.Range("F6").Value = 246
$ws4.Range("F8").Value = 7367
$ws4.Range("F10").Value = 67
$ws4.Range("F12").Value = 3154
$ws4.Range("F14").Value = 577
$ws4.Range("F15").Value = 587
$ws4.Range("F20").Value = 304
$ws4.Range("F23").Value = 746
$ws4.Range("F26").Value = 198
$ws4.Range("F33").Value = 362
$ws4.Range("F35").Value = 1069
$ws4.Range("F37").Value = 115
$ws4.Range("F38").Value = 2101
$ws4.Range("F39").Value = 616
$ws4.Range("F40").Value = 24
